# The presentation currently uses the "Integral" (Red Violet) design/theme
# for its slides (ppt/theme/theme2.xml, wired to the slide master), while a
# second, otherwise-unused theme part (ppt/theme/theme1.xml, wired only to
# the notes master) still carries the stock "Office Theme" colour scheme.
#
# The authored edit swaps which colour scheme is "live": the presentation's
# design is changed over to the plain "Office Theme" palette. We reproduce
# that by rewriting every slot of the active theme's 12-colour scheme
# (ThemeColorScheme, reachable off any Slide) from the old "Red Violet"
# values to the standard "Office" values.
#
# PowerPoint's ThemeColor.RGB is an OLE_COLOR (0x00BBGGRR), i.e. byte-order
# reversed from the usual "RRGGBB" hex notation, so every literal below is
# R + G*256 + B*65536 for the target "RRGGBB" colour.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# 1 = dk1      000000
$tcs.Item(1).RGB  = 0
# 2 = lt1      FFFFFF
$tcs.Item(2).RGB  = 16777215
# 3 = dk2      44546A
$tcs.Item(3).RGB  = 6968388
# 4 = lt2      E7E6E6
$tcs.Item(4).RGB  = 15132391
# 5 = accent1  5B9BD5
$tcs.Item(5).RGB  = 13998939
# 6 = accent2  ED7D31
$tcs.Item(6).RGB  = 3243501
# 7 = accent3  A5A5A5
$tcs.Item(7).RGB  = 10855845
# 8 = accent4  FFC000
$tcs.Item(8).RGB  = 49407
# 9 = accent5  4472C4
$tcs.Item(9).RGB  = 12874308
# 10 = accent6 70AD47
$tcs.Item(10).RGB = 4697456
# 11 = hlink   0563C1
$tcs.Item(11).RGB = 12673797
# 12 = folHlink 954F72
$tcs.Item(12).RGB = 7491477
